$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of data at row 30 (next empty row after existing data)
$row = 30
$ws.Cells.Item($row, 1).Value = "Globo"
$ws.Cells.Item($row, 2).Value = "Inter TV Rural"
$ws.Cells.Item($row, 3).Value = "Agricultura"
$ws.Cells.Item($row, 4).Value = "2025-04-03T10:21"
$ws.Cells.Item($row, 5).Value = "Positivo"
$ws.Cells.Item($row, 6).Value = "testeeeeeeeeee"
